# Oppgave 1, 2 og 3 fungerer optimalt
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "foresatt": add two guardian rows (row 2 and row 3)
# ---------------------------------------------------------------
$wsForesatt = $wb.Worksheets.Item("foresatt")

$wsForesatt.Cells.Item(2, 1).Value = 0
$wsForesatt.Cells.Item(2, 2).Value = 2
$wsForesatt.Cells.Item(2, 3).Value = "Solveig Imsdal"
$wsForesatt.Cells.Item(2, 4).Value = "Bekkeveien 100"
$wsForesatt.Cells.Item(2, 5).NumberFormat = "@"
$wsForesatt.Cells.Item(2, 5).Value = "91997087"
$wsForesatt.Cells.Item(2, 6).NumberFormat = "@"
$wsForesatt.Cells.Item(2, 6).Value = "09079233221"

$wsForesatt.Cells.Item(3, 1).Value = 1
$wsForesatt.Cells.Item(3, 2).Value = 1
$wsForesatt.Cells.Item(3, 3).Value = "Ole Nordmann"
$wsForesatt.Cells.Item(3, 4).Value = "Alvestien 39"
$wsForesatt.Cells.Item(3, 5).NumberFormat = "@"
$wsForesatt.Cells.Item(3, 5).Value = "91997087"
$wsForesatt.Cells.Item(3, 6).NumberFormat = "@"
$wsForesatt.Cells.Item(3, 6).Value = "09079089332"

# column A carries the same bold/bordered "index" look as the header row
$wsForesatt.Range("B1").Copy()
$wsForesatt.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Sheet "barnehage": drop the stray "index" column (column B)
# ---------------------------------------------------------------
$wsBarnehage = $wb.Worksheets.Item("barnehage")
$wsBarnehage.Columns.Item(2).Delete()

# ---------------------------------------------------------------
# Sheet "barn": add one child row (row 2)
# ---------------------------------------------------------------
$wsBarn = $wb.Worksheets.Item("barn")
$wsBarn.Cells.Item(2, 1).Value = 0
$wsBarn.Cells.Item(2, 2).Value = 1
$wsBarn.Cells.Item(2, 3).NumberFormat = "@"
$wsBarn.Cells.Item(2, 3).Value = "09012356472"
$wsBarn.Range("B1").Copy()
$wsBarn.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Sheet "soknad": update existing rows + append a new application row
# ---------------------------------------------------------------
$wsSoknad = $wb.Worksheets.Item("soknad")

# Row 2
$wsSoknad.Cells.Item(2, 2).Value = 4
$wsSoknad.Cells.Item(2, 3).Value = 1
$wsSoknad.Cells.Item(2, 12).NumberFormat = "@"
$wsSoknad.Cells.Item(2, 12).Value = "2020-12-09"
$wsSoknad.Cells.Item(2, 13).NumberFormat = "@"
$wsSoknad.Cells.Item(2, 13).Value = "1100000"

# Row 3
$wsSoknad.Cells.Item(3, 2).Value = 3
$wsSoknad.Cells.Item(3, 3).Value = 1
$wsSoknad.Cells.Item(3, 6).Value = "on"
$wsSoknad.Cells.Item(3, 7).Value = "on"
$wsSoknad.Cells.Item(3, 8).Value = "on"
$wsSoknad.Cells.Item(3, 9).Value = "hei"
$wsSoknad.Cells.Item(3, 11).Value = "on"
$wsSoknad.Cells.Item(3, 12).NumberFormat = "@"
$wsSoknad.Cells.Item(3, 12).Value = "2025-05-07"
$wsSoknad.Cells.Item(3, 13).Value = 1100000

# Row 4
$wsSoknad.Cells.Item(4, 2).Value = 2
$wsSoknad.Cells.Item(4, 3).Value = 1
$wsSoknad.Cells.Item(4, 6).Value = "on"
$wsSoknad.Cells.Item(4, 7).Value = "on"
$wsSoknad.Cells.Item(4, 8).Value = "on"
$wsSoknad.Cells.Item(4, 9).Value = "hei"
$wsSoknad.Cells.Item(4, 11).Value = "on"
$wsSoknad.Cells.Item(4, 12).NumberFormat = "@"
$wsSoknad.Cells.Item(4, 12).Value = "2025-05-07"
$wsSoknad.Cells.Item(4, 13).Value = 1100000

# Row 5 (new application)
$wsSoknad.Cells.Item(5, 1).Value = 3
$wsSoknad.Cells.Item(5, 2).Value = 1
$wsSoknad.Cells.Item(5, 3).Value = 1
$wsSoknad.Cells.Item(5, 4).Value = 2
$wsSoknad.Cells.Item(5, 5).Value = 1
$wsSoknad.Cells.Item(5, 6).Value = "on"
$wsSoknad.Cells.Item(5, 7).Value = "on"
$wsSoknad.Cells.Item(5, 8).Value = "on"
$wsSoknad.Cells.Item(5, 9).Value = "hei"
$wsSoknad.Cells.Item(5, 10).Value = "Sunshine Preschool"
$wsSoknad.Cells.Item(5, 11).Value = "on"
$wsSoknad.Cells.Item(5, 12).NumberFormat = "@"
$wsSoknad.Cells.Item(5, 12).Value = "2025-05-07"
$wsSoknad.Cells.Item(5, 13).Value = 1100000

$wsSoknad.Range("A4").Copy()
$wsSoknad.Range("A5").PasteSpecial(-4122)
